$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Regenerated "K" (strike count) column values after switching save_data
# to compute K from the data instead of the legacy Strike# field.
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("G8").Value = 1
$ws.Range("G9").Value = 2
$ws.Range("G10").Value = 2
$ws.Range("G12").Value = 0
$ws.Range("G13").Value = 1
$ws.Range("G14").Value = 1
$ws.Range("G15").Value = 1
